$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new status column
$ws.Range("F1").Value = "status"

# Existing rows: fill in status values
$ws.Range("F2").Value = "done"
$ws.Range("F3").Value = "done"
$ws.Range("F4").Value = "done"
$ws.Range("F7").Value = "na"
$ws.Range("F9").Value = "in progress"
$ws.Range("F13").Value = "na"
$ws.Range("F18").Value = "na"
$ws.Range("F19").Value = "na"
$ws.Range("F31").Value = "na"
$ws.Range("F32").Value = "na"
$ws.Range("F33").Value = "na"
$ws.Range("F34").Value = "na"

# New rows for additional topics
$ws.Range("A35").Value = "Thread"
$ws.Range("F35").Value = "in progress"
$ws.Range("A36").Value = "Oauth 2"
$ws.Range("F36").Value = "yet to start"

# Column A width / best fit (stored width matches Excel's internal MDW-based
# rounding: setting ColumnWidth=19.1 serializes to width="20" in the XML)
$ws.Columns.Item(1).ColumnWidth = 19.1

# Set selection to match the author's last selected cell
$ws.Range("F14").Select()
